# Weekly update: insert 3 new price rows (one per "Calidad": Especial, Primera,
# Segunda) at the top of the data block (row 865), pushing the existing
# historical rows down by three (865-924 -> 868-927).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before the current row 865, shifting everything
# below (including formatting) down by three rows.
$ws.Range("A865:A867").EntireRow.Insert()

# Common (unchanged) column values shared by the three new rows.
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100101
$producto = "Berries"
$categoriaId = 100112025
$categoria = "Frutilla"
$variedad = "Sin especificar"
$unidad = "`$/bandeja 7 kilos"
$origen = "Provincia de Melipilla"
$kgUnidad = 7

function Set-NewRow($row, $calidad, $fecha, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioMax
    $ws.Cells.Item($row, 16).Value = $precioProm
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

Set-NewRow 865 "Especial" 44610 1500 6000 7000 6500 929
Set-NewRow 866 "Primera"  44610 1850 4000 5000 4500 643
Set-NewRow 867 "Segunda"  44610 450  2500 3500 3000 429
